$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "cUDuz984"
$ws.Range("B2").Value = 23102882
$ws.Range("C2").Value = "vjauicx43"
$ws.Range("D2").Value = "jCh#5%8T"
$ws.Range("F2").Value = "QzrdpkPC"
$ws.Range("G2").Value = "JAdc"

# Row 3
$ws.Range("A3").Value = "oaFQH721"
$ws.Range("B3").Value = 23102881
$ws.Range("C3").Value = "uhdtmtz68"
$ws.Range("D3").Value = "F#8Cn5b$"
$ws.Range("F3").Value = "ayCLhEVp"
$ws.Range("G3").Value = "zKgE"
